$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.878.56"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "2.599.13"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'591.85"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'151.37"
$ws.Range("E6").Value = "  -3.38%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "2.597.11"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("E10").Value = "  -3.37%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -2.80%  "
$ws.Range("E13").Value = "  -4.53%  "
$ws.Range("D14").Value = "'27.27"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "3.074.12"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("E16").Value = "  -5.13%  "
$ws.Range("D17").Value = "66.743.11"
$ws.Range("E17").Value = "  -1.24%  "
$ws.Range("D18").Value = "2.624.29"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'362.30"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("E21").Value = "  -5.88%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'2.03"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'9.88"
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'66.18"
$ws.Range("E26").Value = "  -3.05%  "
$ws.Range("D27").Value = "2.736.25"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'576.01"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("D30").Value = "0.0₃0999"
$ws.Range("E30").Value = "  -5.14%  "
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("E32").Value = "  -4.97%  "
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -8.67%  "
$ws.Range("E36").Value = "  -5.29%  "
$ws.Range("E37").Value = "  -3.42%  "
$ws.Range("D38").Value = "'155.80"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("E40").Value = "  -2.47%  "
$ws.Range("D41").Value = "'5.19"
$ws.Range("E41").Value = "  -4.88%  "
$ws.Range("E42").Value = "  -5.19%  "
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'40.78"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'16.41"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "'153.74"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "0.0₆0284"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("D50").Value = "'21.23"
$ws.Range("E50").Value = "  +0.49%  "
$ws.Range("D51").Value = "'0.612"
$ws.Range("E51").Value = "  -3.83%  "
